$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 85, shifting existing rows 85-92 down to 87-94.
$ws.Rows("85:86").Insert()

# --- New row 85 (weekly update, new record) ---
$ws.Range("A85").Value = 6
$ws.Range("B85").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C85").Value = "Metropolitana"
$ws.Range("D85").Value = 44918
$ws.Range("E85").Value = 13
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100101
$ws.Range("H85").Value = "Berries"
$ws.Range("I85").Value = 100101008
$ws.Range("J85").Value = "Mora"
$ws.Range("K85").Value = "Sin especificar"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 200
$ws.Range("N85").Value = 4000
$ws.Range("O85").Value = 4000
$ws.Range("P85").Value = 4000
$ws.Range("Q85").Value = "`$/bandeja 2 kilos"
$ws.Range("R85").Value = "Provincia de Curicó"
$ws.Range("S85").Value = 2000
$ws.Range("T85").Value = 2

# --- New row 86 (weekly update, new record) ---
$ws.Range("A86").Value = 6
$ws.Range("B86").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44918
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100101008
$ws.Range("J86").Value = "Mora"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Segunda"
$ws.Range("M86").Value = 150
$ws.Range("N86").Value = 3000
$ws.Range("O86").Value = 3000
$ws.Range("P86").Value = 3000
$ws.Range("Q86").Value = "`$/bandeja 2 kilos"
$ws.Range("R86").Value = "Provincia de Curicó"
$ws.Range("S86").Value = 1500
$ws.Range("T86").Value = 2
